$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.505614041169197
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 17.36656647638025

# Row 3
$ws.Range("B3").Value = 0.006876353814593728
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 71535.21976971527

# Row 4
$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 86.29678392075563
$ws.Range("D4").Value = 2938.103010863317
$ws.Range("E4").Value = 1594453305621061000
$ws.Range("G4").Value = 1594453305621064000
